$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Create the three new character styles referenced by the new runs.
# ---------------------------------------------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# ---------------------------------------------------------------------
# 2. Apply the relevant character style to every run that carries the
#    matching text (the document repeats the "Dates de la campanya..."
#    paragraph four times, and the other two paragraphs once each).
# ---------------------------------------------------------------------

function Apply-StyleToAllMatches($searchText, $styleName) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $searchText
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 0
    $rng.Find.MatchWildcards = $false
    while ($rng.Find.Execute()) {
        $rng.Style = $styleName
        $rng.Collapse(0)
    }
}

$datesText = "Dates de la campanya 2022 en qu" + [char]232 + " usem la  Constel" + [char]183 + "laci" + [char]243 + " de Peg" + [char]224 + "s 8-17 d'octubre, 7-16 de novembre,"
Apply-StyleToAllMatches $datesText "GaNStyle"

$esteuText = "Esteu participant en una campanya mundial per observar i anotar la brillantor de les estrelles m" + [char]233 + "s febles que es poden veure, com a mitj" + [char]224 + " per mesurar la contaminaci" + [char]243 + " lum" + [char]237 + "nica en un lloc determinat. Localitzant i observant la  Constel" + [char]183 + "laci" + [char]243 + " de Peg" + [char]224 + "s a la nit i comparant la brillantor de les estrelles del cel amb la brillantor que indiquen els mapes, gent de tot el m" + [char]243 + "n aprendran com els llums de la seva zona contribueixen a augmentar la contaminaci" + [char]243 + " lum" + [char]237 + "nica. Les vostres aportacions a la base de dades activa faran palesa la visibilitat del cel nocturn."
Apply-StyleToAllMatches $esteuText "GaNParagraph"

$jenikText = "Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
Apply-StyleToAllMatches $jenikText "GaNLinks"
